$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-10-28"

$ws.Range("B1").Value = "October 2021 (through October 28)"

$ws.Range("V2").Value = 8
$ws.Range("AF5").Value = 1
$ws.Range("L8").Value = 6
$ws.Range("AP9").Value = 3
$ws.Range("V10").Value = 2
$ws.Range("AP10").Value = 5
$ws.Range("AZ10").Value = 4
$ws.Range("AZ13").Value = 4
$ws.Range("AP15").Value = 2
$ws.Range("BJ21").Value = 1
$ws.Range("V25").Value = 1
$ws.Range("V33").Value = 2
$ws.Range("B36").Value = 3
$ws.Range("AP38").Value = 5
$ws.Range("L39").Value = 2
$ws.Range("AP39").Value = 1
$ws.Range("AZ41").Value = 1
$ws.Range("AP47").Value = 2
$ws.Range("AF49").Value = 1
$ws.Range("L64").Value = 2
$ws.Range("AP72").Value = 1
$ws.Range("AP93").Value = 1
